# Update on cost estimation
$wb = $excel.ActiveWorkbook

$costSheet = $wb.Worksheets.Item("Cost")
$heightSheet = $wb.Worksheets.Item("Height")

# --- Cost sheet -----------------------------------------------------------
# The "Camera 8MP V2" line is removed; every row below it shifts up by one
# (item index numbers in columns A/E are untouched), the "Touch Screen" row
# gets new cost figures, and the old "Misc" line is replaced by a new
# "Infrared Sensor" line. The totals formulas move up a row and now sum
# C2:C9 / G2:G9 instead of C2:C10 / G2:G10.

$costSheet.Range("B5").Value = "DHT11"
$costSheet.Range("C5").Value = 9.9
$costSheet.Range("F5").Value = "DHT11"
$costSheet.Range("G5").Value = 1.58

$costSheet.Range("B6").Value = "Rubber Hose"
$costSheet.Range("C6").Value = 5
$costSheet.Range("F6").Value = "Rubber Hose"
$costSheet.Range("G6").Value = 5

$costSheet.Range("B7").Value = "Plastic Container"
$costSheet.Range("C7").Value = 4.9400000000000004
$costSheet.Range("F7").Value = "Plastic Container"
$costSheet.Range("G7").Value = 4.9400000000000004

$costSheet.Range("B8").Value = "Touch Screen"
$costSheet.Range("C8").Value = 63
$costSheet.Range("F8").Value = "Touch Screen"
$costSheet.Range("G8").Value = 42

$costSheet.Range("B9").Value = "Infrared Sensor"
$costSheet.Range("C9").Value = 18.95
$costSheet.Range("F9").Value = "Infrared Sensor"
$costSheet.Range("G9").Value = 18.95

# Old row 10 (index 9, "Misc") and the old totals row 11 collapse into the
# new totals row 10.
$costSheet.Range("A10").Clear()
$costSheet.Range("B10").Clear()
$costSheet.Range("E10").Clear()
$costSheet.Range("F10").Clear()
$costSheet.Range("C10").Formula = "=SUM(C2:C9)"
$costSheet.Range("G10").Formula = "=SUM(G2:G9)"

$costSheet.Range("A11:G11").Clear()

# --- View state -------------------------------------------------------
# "Cost" becomes the active/selected tab (was "Height"), zoomed to 90%,
# with G9 selected; "Height" keeps D14 selected for when it regains focus.
$heightSheet.Activate()
$heightSheet.Range("D14").Select() | Out-Null

$costSheet.Activate()
$excel.ActiveWindow.Zoom = 90
$costSheet.Range("G9").Select() | Out-Null

$wb.Save()
